# Insert a new data row at row 84 (pushing existing rows 84-92 down to 85-93)
# and populate it with the new weekly data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 84, shifting rows 84:92 down to 85:93
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with values, matching the style/format
# of the other data rows (column D uses style index "2" in the original file,
# i.e. the same date number format as the rows above/below it).
$ws.Range("A84").Value = 8
$ws.Range("B84").Value = "Terminal La Palmera de La Serena"
$ws.Range("C84").Value = "Coquimbo"
$ws.Range("D84").Value = 44504
$ws.Range("D84").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E84").Value = 4
$ws.Range("F84").Value = 100112044
$ws.Range("G84").Value = "Perejil"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 3200
$ws.Range("K84").Value = 1300
$ws.Range("L84").Value = 1500
$ws.Range("M84").Value = 1400
$ws.Range("N84").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O84").Value = "Provincia del Elquí"
$ws.Range("P84").Value = 933
$ws.Range("Q84").Value = 1.5
$ws.Range("R84").Value = "Hortaliza"

$wb.Save()
